$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.865.74"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "3.137.48"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.58"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "626.58"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  +7.17%  "
$ws.Range("E8").Value = "  +5.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "3.137.62"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.760"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.205"
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.82"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "91.348.51"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "3.715.53"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "3.090.64"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.80"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.71"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000217"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.91"
$ws.Range("E22").Value = "  +3.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "453.10"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.94"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.16"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.179"
$ws.Range("E30").Value = "  +11.69%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.233"
$ws.Range("E31").Value = "  +14.82%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.124"
$ws.Range("E32").Value = "  +43.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.19"
$ws.Range("E33").Value = "  -6.77%  "
$ws.Range("E34").Value = "  +26.16%  "
$ws.Range("E35").Value = "  +8.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.92"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.54"
$ws.Range("E37").Value = "  +6.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.20"
$ws.Range("E38").Value = "  +26.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "503.19"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.425"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.22"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.58"
$ws.Range("E47").Value = "  +5.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.702"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.58"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.16"
$ws.Range("E51").Value = "  +0.18%  "
